$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 11, pushing the existing data (old rows 11-56) down to rows 13-58.
$ws.Rows("11:12").Insert()

# Fill in the new weekly data row for "Primera" quality.
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C11").Value = "Arica y Parinacota"
$ws.Range("D11").Value = 44592
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = 100112021
$ws.Range("G11").Value = "Ají"
$ws.Range("H11").Value = "Inferno"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 140
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 14500
$ws.Range("N11").Value = "$/caja 15 kilos"
$ws.Range("O11").Value = "Región de Arica y Parinacota"
$ws.Range("P11").Value = 967
$ws.Range("Q11").Value = 15
$ws.Range("R11").Value = "Hortaliza"

# Fill in the new weekly data row for "Segunda" quality.
$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C12").Value = "Arica y Parinacota"
$ws.Range("D12").Value = 44592
$ws.Range("E12").Value = 15
$ws.Range("F12").Value = 100112021
$ws.Range("G12").Value = "Ají"
$ws.Range("H12").Value = "Inferno"
$ws.Range("I12").Value = "Segunda"
$ws.Range("J12").Value = 130
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 11000
$ws.Range("M12").Value = 10500
$ws.Range("N12").Value = "$/caja 15 kilos"
$ws.Range("O12").Value = "Región de Arica y Parinacota"
$ws.Range("P12").Value = 700
$ws.Range("Q12").Value = 15
$ws.Range("R12").Value = "Hortaliza"
